# ModMatrixTable.xlsx edit: rename Sheet1 -> Main, add a new "Comments" sheet,
# and restructure the data table on Main (add Mix columns / Morph rows).

$wb = $excel.ActiveWorkbook
$main = $wb.ActiveSheet

# --- Rename the original sheet and wipe its contents/formatting so we can
# --- rebuild the table cleanly in its new layout. ---
$main.Name = "Main"
$main.UsedRange.Clear() | Out-Null

# --- Row 1: header labels ---
$main.Range("A1").Value = "Scale.Coeff"
$main.Range("B1").Value = "MixX"
$main.Range("C1").Value = "MixY"
$main.Range("D1").Value = "Freq"
$main.Range("E1").Value = "Cutoff"
$main.Range("F1").Value = "LFOFreq"

# --- Row 2: MorphX ---
$main.Range("A2").Value = "MorphX"
$main.Range("B2").Value = 1
$main.Range("C2").Value = 0
$main.Range("D2").Value = 0
$main.Range("E2").Value = 0
$main.Range("F2").Value = 0

# --- Row 3: MorphY ---
$main.Range("A3").Value = "MorphY"
$main.Range("B3").Value = 0
$main.Range("C3").Value = 1
$main.Range("D3").Value = 0
$main.Range("E3").Value = 0
$main.Range("F3").Value = 0

# --- Row 4: Expression ---
$main.Range("A4").Value = "Expression"
$main.Range("B4").Value = 0
$main.Range("C4").Value = 0
$main.Range("D4").Value = 100
$main.Range("E4").Value = 500
$main.Range("F4").Value = 5

# --- Row 5: LFO ---
$main.Range("A5").Value = "LFO"
$main.Range("B5").Value = 0
$main.Range("C5").Value = 0
$main.Range("D5").Value = 10
$main.Range("E5").Value = 50
$main.Range("F5").Value = 0

# --- Bold styling on header row, label column, and the assorted blank
# --- formatting-only cells that were carried over from the source sheet. ---
$main.Range("A1:F1").Font.Bold = $true
$main.Range("G1").Font.Bold = $true
$main.Range("L1:N1").Font.Bold = $true
$main.Range("O1").Font.Bold = $true
$main.Range("Q1:T1").Font.Bold = $true

$main.Range("A2").Font.Bold = $true
$main.Range("G2").Font.Bold = $true
$main.Range("L2").Font.Bold = $true
$main.Range("Q2").Font.Bold = $true

$main.Range("A3").Font.Bold = $true
$main.Range("L3").Font.Bold = $true
$main.Range("Q3").Font.Bold = $true

$main.Range("A4").Font.Bold = $true

$main.Range("A5").Font.Bold = $true
$main.Range("G5:J5").Font.Bold = $true
$main.Range("L5:O5").Font.Bold = $true
$main.Range("Q5:T5").Font.Bold = $true

$main.Range("G6").Font.Bold = $true
$main.Range("L6").Font.Bold = $true
$main.Range("Q6").Font.Bold = $true

$main.Range("A7:B7").Font.Bold = $true
$main.Range("G7").Font.Bold = $true
$main.Range("L7").Font.Bold = $true
$main.Range("Q7").Font.Bold = $true

# --- Column widths ---
$main.Columns.Item(1).ColumnWidth = 17.166666666666668
$main.Columns.Item(2).ColumnWidth = 17.166666666666668
$main.Columns.Item(9).ColumnWidth = 11.666666666666666
$main.Columns.Item(10).ColumnWidth = 11.666666666666666
$main.Columns.Item(12).ColumnWidth = 13.5

# --- Add the new "Comments" sheet right after Main ---
$comments = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $main)
$comments.Name = "Comments"

$comments.Range("D1").Value = "Noisy"
$comments.Range("D1").Font.Bold = $true
$comments.Range("E1").Value = "Roomy"

$comments.Range("D2").Value = "1 - morphX"
$comments.Range("E2").Value = """=morphX"""

$comments.Range("A4").Value = "x"
$comments.Range("B4").Value = "0-1"
$comments.Range("A5").Value = "y"
$comments.Range("B5").Value = "0-1"

$comments.Range("D26").Select() | Out-Null

# --- Re-activate Main, restore zoom and selection ---
$main.Activate() | Out-Null
$main.Application.ActiveWindow.Zoom = 160
$main.Range("A9").Select() | Out-Null
